$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BP-813: Affiliate Mapping for True Independent Stations
# Rename columns: "Affiliation Mismatch Note" -> "IsTrueIND" and "SalesGroupName" -> "RepFirm"
# Also extend the header/body formatting (thin border) across columns I:J to match the rest of the table.

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("I2:J5").PasteSpecial(-4122)

$ws.Range("H1").Value = "IsTrueIND"
$ws.Range("J1").Value = "RepFirm"
